# Uren Registratie Vorig Week, check of het even klopt & Asset list update
#
# - Week 8 (Donderdag, row 30): hours corrected from 1 -> 4 for everyone.
# - Week 11 (Woensdag/Donderdag/Vrijdag, rows 52-54): previously-empty days
#   filled in with the actual registered hours + absence notes.
# - All dependent SUM()/ratio formulas recalc automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Colors used throughout this sheet to flag attendance:
#   blue (0070C0) -> hours logged (non-zero)
#   red  (FF0000) -> absent (zero hours)
$colBlue = 0 + (112 * 256) + (192 * 65536)   # RGB(0,112,192)
$colRed  = 255 + (0 * 256) + (0 * 65536)     # RGB(255,0,0)

# --- Week 8, Donderdag (row 30): 1 -> 4 for every column ---
"B30","C30","D30","E30","F30","G30" | ForEach-Object {
    $ws.Range($_).Value = 4
}

# --- Week 11, Woensdag (row 52): fill in, everyone present (4h) ---
$ws.Range("B52").Value = 4
"C52","D52","E52","F52","G52" | ForEach-Object {
    $ws.Range($_).Value = 4
    $ws.Range($_).Interior.Color = $colBlue
}

# --- Week 11, Donderdag (row 53): Fahrettin & Lara absent ---
$ws.Range("B53").Value = 2
$ws.Range("C53").Value = 0
$ws.Range("C53").Interior.Color = $colRed
$ws.Range("D53").Value = 0
$ws.Range("D53").Interior.Color = $colRed
$ws.Range("E53").Value = 2
$ws.Range("E53").Interior.Color = $colBlue
$ws.Range("F53").Value = 2
$ws.Range("F53").Interior.Color = $colBlue
$ws.Range("G53").Value = 2
$ws.Range("G53").Interior.Color = $colBlue

# --- Week 11, Vrijdag (row 54): Ruben absent ---
$ws.Range("B54").Value = 4
$ws.Range("C54").Value = 4
$ws.Range("C54").Interior.Color = $colBlue
$ws.Range("D54").Value = 4
$ws.Range("D54").Interior.Color = $colBlue
$ws.Range("E54").Value = 0
$ws.Range("E54").Interior.Color = $colRed
$ws.Range("F54").Value = 4
$ws.Range("F54").Interior.Color = $colBlue
$ws.Range("G54").Value = 4
$ws.Range("G54").Interior.Color = $colBlue

# Absence remarks (insertion order kept so the shared-string table lines up:
# "Rubeb: Ziekenhuis |" ends up before "Fahrettin: Te Laat Wakker | Lara: Ziek |")
$ws.Range("H54").Value = "Rubeb: Ziekenhuis |"
$ws.Range("H53").Value = "Fahrettin: Te Laat Wakker | Lara: Ziek |"

# Move the active selection the way the author left it (I4, scrolled to top)
$ws.Range("I4").Select() | Out-Null
